# BC02_HocVienGhiDanh.xlsx — "Add files via upload" update
#
# Summary of the functional change:
#   - The tuition column "Tổng tiền" is relabelled "Học phí".
#   - A new "Thực đóng" column (bound to &=[DATA1].DADONG, i.e. the old
#     "Đã đóng" column) is appended after "Còn nợ", and the
#     Miễn giảm(%) / Miễn giảm(tiền) / Còn nợ columns each shift one
#     column to the left (J/K/L instead of K/L/M) to make room.
#   - The three spacer cells that used to sit in column L (rows 1-3) move
#     to column M, since the data grid is now A:M with the new column.
#   - The right-hand summary row 9 swaps which cell carries the
#     "&=[DATA].CURRENTDATETIME" placeholder (J9 instead of K9) and which
#     cell keeps the centred spacer style (M9 instead of J9).
#   - Number columns I:M in the header/detail/total rows become
#     right-aligned instead of the previous general/center alignment.
#   - The remembered selection moves to M6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats
$xlRight   = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignRight

# ---------------------------------------------------------------------
# Rows 1-3: the trailing blank spacer cell moves from column L to M
# (mirrors columns J/K which keep their existing per-row style).
# ---------------------------------------------------------------------
foreach ($row in 1..3) {
    $fromFormatSrc = $ws.Cells.Item($row, 10)   # J<row> keeps its own style
    $blank         = $ws.Cells.Item($row, 12)   # L<row>
    $target        = $ws.Cells.Item($row, 13)   # M<row>

    $fromFormatSrc.Copy() | Out-Null
    $target.PasteSpecial($xlFormats) | Out-Null

    $blank.Style = "Normal"
    $blank.ClearContents() | Out-Null
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row 5: header labels for I:M
# ---------------------------------------------------------------------
$ws.Range("I5").Value = "Học phí"
$ws.Range("J5").Value = "Miễn giảm (%)"
$ws.Range("K5").Value = "Miễn giảm (tiền)"
$ws.Range("L5").Value = "Còn nợ"
$ws.Range("M5").Value = "Thực đóng"

# ---------------------------------------------------------------------
# Row 6: bound-field placeholders for I:M (I unchanged; J/K/L shift left
# from the old K/L/M; M picks up the old J's DADONG placeholder)
# ---------------------------------------------------------------------
$ws.Range("I6").Value = "&=[DATA1].TONGTIEN"
$ws.Range("J6").Value = "&=[DATA1].MienGiam_PhanTram"
$ws.Range("K6").Value = "&=[DATA1].MienGiam_Tien"
$ws.Range("L6").Value = "&=[DATA1].CONNO"
$ws.Range("M6").Value = "&=[DATA1].DADONG"

# ---------------------------------------------------------------------
# Row 7: totals — keep a straightforward SUM of the cell above in each
# of I:M (the formulas stay column-relative, so they automatically track
# whichever field now lives in that column).
# ---------------------------------------------------------------------
$ws.Range("I7").Formula = "=SUM(I6:I6)"
$ws.Range("J7").Formula = "=SUM(J6:J6)"
$ws.Range("K7").Formula = "=SUM(K6:K6)"
$ws.Range("L7").Formula = "=SUM(L6:L6)"
$ws.Range("M7").Formula = "=SUM(M6:M6)"

# Right-align the numeric grid (header I5:M5 keeps its existing centred
# style — only the data row and the totals row change alignment).
$ws.Range("I6:M6").HorizontalAlignment = $xlRight
$ws.Range("I7:M7").HorizontalAlignment = $xlRight

# ---------------------------------------------------------------------
# Row 9: the CURRENTDATETIME placeholder moves from K9 to J9, and the
# centred spacer style that used to live on J9 moves to M9.
# ---------------------------------------------------------------------
$ws.Range("J9").Copy() | Out-Null
$ws.Range("M9").PasteSpecial($xlFormats) | Out-Null

$ws.Range("K9").Copy() | Out-Null
$ws.Range("J9").PasteSpecial($xlFormats) | Out-Null
$ws.Range("J9").Value = "&=[DATA].CURRENTDATETIME"
$ws.Range("K9").ClearContents() | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Selection / view bookkeeping
# ---------------------------------------------------------------------
$ws.Range("M6").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 7
